# Adds a new worksheet "Planilha2" with UST consolidation data (estimativas e
# consolidado UST) after "Planilha1", and makes it the active sheet.
# Commit: "inclui arquivos estimativas e consolidado UST"

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet, then move it so it is placed after Planilha1 ---
$ws2 = $wb.Worksheets.Add()
$ws2.Name = "Planilha2"
$ws1 = $wb.Worksheets.Item("Planilha1")
$ws2.Move($null, $ws1)

# NOTE: worksheet object variables in this host are resolved by position
# (index), and both Add() and Move() re-index the sheet collection under the
# hood. Always re-fetch fresh references by name right after such calls,
# before doing anything else with them.
$ws1 = $wb.Worksheets.Item("Planilha1")
$ws2 = $wb.Worksheets.Item("Planilha2")

# --- Header row ---
$ws2.Range("B1").Value = "historias"
$ws2.Range("C1").Value = "sprints"
$ws2.Range("D1").Value = "semanas"
$ws2.Range("E1").Value = "equipe_num"
$ws2.Range("F1").Value = "usts_equipe"
$ws2.Range("G1").Value = "usts_individuais"
$ws2.Range("H1").Value = "usts_totais"

# --- Row 2: Gabriel ---
$ws2.Range("A2").Value = "Gabriel"
$ws2.Range("B2").Value = 10
$ws2.Range("C2").Value = 5
$ws2.Range("D2").Value = 5
$ws2.Range("E2").Value = 4
$ws2.Range("F2").Value = 640
$ws2.Range("G2").Value = 2260
$ws2.Range("H2").Value = 2900

# --- Row 3: Francisco ---
$ws2.Range("A3").Value = "Francisco"
$ws2.Range("C3").Value = 6
$ws2.Range("D3").Value = 12
$ws2.Range("E3").Value = 5
$ws2.Range("F3").Value = 162
$ws2.Range("G3").Formula = "=H3-F3"
$ws2.Range("H3").Value = 850

# --- Row 4: Maurício ---
$ws2.Range("A4").Value = "Maurício"

# --- Bold styling for name column / header-like cells (matches style index 1
#     used on Planilha1's header row) ---
$ws2.Range("A2:A4").Font.Bold = $true
$ws2.Range("B2:B4").Font.Bold = $true

# --- J/K side-note columns (estimates breakdown per person) ---
$ws2.Range("J3").Value = "Francisco"
$ws2.Range("K3").Value = "Gabriel"

$ws2.Range("J4").Value = "documentação e apresentação (80), "
$ws2.Range("K4").Value = "documentação e apresentação (48), "

$ws2.Range("J5").Value = "identidade visual (30), "
$ws2.Range("K5").Value = "implementação front-end ou back end (19), "

$ws2.Range("J6").Value = "acompanhamento do scrum master (28), "
$ws2.Range("K6").Value = "identidade visual (15), "

$ws2.Range("J7").Value = "arquitetura/definições(20), "
$ws2.Range("K7").Value = "entendimento, refinamento, escrita e validação (14), "

$ws2.Range("J8").Value = "implementação front-end (19), "
$ws2.Range("K8").Value = "prototipação (9), "

$ws2.Range("J9").Value = "entendimento, refinamento, escrita e validação (14), "
$ws2.Range("K9").Value = "preparação e implementação (8), "

$ws2.Range("J10").Value = "guia usabilidade (10),"
$ws2.Range("K10").Value = "guia usabilidade (5)"

$ws2.Range("J11").Value = "prototipação (9),"

$ws2.Range("J12").Value = "preparação e implementação (8)"

# --- Column widths (best-fit, matching target <cols>) ---
$ws2.Columns.Item(5).ColumnWidth = 12.28515625
$ws2.Columns.Item(6).ColumnWidth = 11.85546875
$ws2.Columns.Item(7).ColumnWidth = 15.42578125
$ws2.Columns.Item(8).ColumnWidth = 10.5703125
$ws2.Columns.Item(10).ColumnWidth = 49.28515625

# --- Selections: this host persists a sheet's <selection> only by replaying
#     the Select() that happened while that sheet was the active one, so set
#     Planilha1's first (while Planilha1 is active), then finish on
#     Planilha2 (left as the active / tabSelected sheet). ---
$ws1.Activate()
$ws1.Range("D1:F1").Select()

$ws2.Activate()
$ws2.Range("B2").Select()
